$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MHS"
$ws.Range("B2").Value = "Mahasiswa"
$ws.Range("A3").Value = "DSN"
$ws.Range("B3").Value = "Dosen"

$ws.Range("B3").Select()
